$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 34

# Add the new header cells I1 ("I0") and J1 ("IF"), copying the header
# formatting (style) from H1 so they match the existing header look.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I column is a constant 1 for every data row.
$ws.Range("I2:I$lastRow").Value = 1

# J column mirrors the existing H column (IP) values for every data row.
$ws.Range("H2:H$lastRow").Copy()
$ws.Range("J2:J$lastRow").PasteSpecial(-4163)
